$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H39").Value = 295.13333
$ws_ALC.Range("I39").Value = 66.40000000000001
$ws_ALC.Range("J39").Value = 409.5
$ws_ALC.Range("K39").Value = 199.2
$ws_ALC.Range("L39").Value = 1228.5
$ws_ALC.Range("M39").Value = 96.79999999999998
$ws_ALC.Range("N39").Value = -1820.5

$ws_ALC.Range("H98").Value = 51224.6
$ws_ALC.Range("I98").Value = 56736.777
$ws_ALC.Range("J98").Value = 1615
$ws_ALC.Range("K98").Value = 56736.777
$ws_ALC.Range("L98").Value = 1615
$ws_ALC.Range("M98").Value = -55238.777
$ws_ALC.Range("N98").Value = -4611

$ws_ALC.Range("H122").Value = 51224.6
$ws_ALC.Range("I122").Value = 56736.777
$ws_ALC.Range("J122").Value = 1615
$ws_ALC.Range("K122").Value = 170210.331
$ws_ALC.Range("L122").Value = 4845
$ws_ALC.Range("M122").Value = -167760.331
$ws_ALC.Range("N122").Value = -9745

$ws_ALC.Range("H127").Value = 1037.4242
$ws_ALC.Range("I127").Value = 670.5
$ws_ALC.Range("J127").Value = 1771.2727
$ws_ALC.Range("K127").Value = 2011.5
$ws_ALC.Range("L127").Value = 5313.8181
$ws_ALC.Range("M127").Value = 2948.5
$ws_ALC.Range("N127").Value = -15233.8181

$ws_ALC.Range("H132").Value = 942620.4399999999
$ws_ALC.Range("I132").Value = 1069058.9
$ws_ALC.Range("J132").Value = 3363.1428
$ws_ALC.Range("K132").Value = 3207176.7
$ws_ALC.Range("L132").Value = 10089.4284
$ws_ALC.Range("M132").Value = -3204646.7
$ws_ALC.Range("N132").Value = -15149.4284

$ws_ALC.Range("H137").Value = 1220.9524
$ws_ALC.Range("I137").Value = 1034.375
$ws_ALC.Range("J137").Value = 1818
$ws_ALC.Range("K137").Value = 3103.125
$ws_ALC.Range("L137").Value = 5454
$ws_ALC.Range("M137").Value = -553.125
$ws_ALC.Range("N137").Value = -10554

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 5521.75
$ws_ARM.Range("I32").Value = 2655.3125
$ws_ARM.Range("K32").Value = 2655.3125
$ws_ARM.Range("M32").Value = -2368.3125

$ws_ARM.Range("H61").Value = 1614.3928
$ws_ARM.Range("I61").Value = 1384.2916
$ws_ARM.Range("J61").Value = 2995
$ws_ARM.Range("K61").Value = 1384.2916
$ws_ARM.Range("L61").Value = 2995
$ws_ARM.Range("M61").Value = -1172.2916
$ws_ARM.Range("N61").Value = -3419

$ws_ARM.Range("H136").Value = 1614.3928
$ws_ARM.Range("I136").Value = 1384.2916
$ws_ARM.Range("J136").Value = 2995
$ws_ARM.Range("K136").Value = 4152.8748
$ws_ARM.Range("L136").Value = 8985
$ws_ARM.Range("M136").Value = -1602.8748
$ws_ARM.Range("N136").Value = -14085

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H62").Value = 60000
$ws_BSM.Range("J62").Value = 60000
$ws_BSM.Range("L62").Value = 60000
$ws_BSM.Range("N62").Value = -61372

$ws_BSM.Range("H64").Value = 625.5
$ws_BSM.Range("I64").Value = 211.2
$ws_BSM.Range("J64").Value = 921.4286
$ws_BSM.Range("K64").Value = 211.2
$ws_BSM.Range("L64").Value = 921.4286
$ws_BSM.Range("M64").Value = 13.80000000000001
$ws_BSM.Range("N64").Value = -1371.4286

$ws_BSM.Range("H65").Value = 60000
$ws_BSM.Range("J65").Value = 60000
$ws_BSM.Range("L65").Value = 180000
$ws_BSM.Range("N65").Value = -186864

$ws_BSM.Range("H67").Value = 625.5
$ws_BSM.Range("I67").Value = 211.2
$ws_BSM.Range("J67").Value = 921.4286
$ws_BSM.Range("K67").Value = 211.2
$ws_BSM.Range("L67").Value = 921.4286
$ws_BSM.Range("M67").Value = 568.8
$ws_BSM.Range("N67").Value = -2481.4286

$ws_BSM.Range("H80").Value = 437.85715
$ws_BSM.Range("I80").Value = 45.25
$ws_BSM.Range("K80").Value = 45.25
$ws_BSM.Range("M80").Value = 952.75

$ws_BSM.Range("H83").Value = 437.85715
$ws_BSM.Range("I83").Value = 45.25
$ws_BSM.Range("K83").Value = 226.25
$ws_BSM.Range("M83").Value = 4765.75

$ws_BSM.Range("H86").Value = 1667.5555
$ws_BSM.Range("I86").Value = 1723.5385
$ws_BSM.Range("J86").Value = 1522
$ws_BSM.Range("K86").Value = 1723.5385
$ws_BSM.Range("L86").Value = 1522
$ws_BSM.Range("M86").Value = -600.5385000000001
$ws_BSM.Range("N86").Value = -3768

$ws_BSM.Range("H89").Value = 1667.5555
$ws_BSM.Range("I89").Value = 1723.5385
$ws_BSM.Range("J89").Value = 1522
$ws_BSM.Range("K89").Value = 8617.692500000001
$ws_BSM.Range("L89").Value = 7610
$ws_BSM.Range("M89").Value = -3001.692500000001
$ws_BSM.Range("N89").Value = -18842

$ws_BSM.Range("H99").Value = 2663.4443
$ws_BSM.Range("I99").Value = 3715
$ws_BSM.Range("J99").Value = 1822.2
$ws_BSM.Range("K99").Value = 3715
$ws_BSM.Range("L99").Value = 1822.2
$ws_BSM.Range("M99").Value = -2217
$ws_BSM.Range("N99").Value = -4818.2

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H16").Value = 4291.273
$ws_CRP.Range("I16").Value = 2741.5715
$ws_CRP.Range("J16").Value = 7003.25
$ws_CRP.Range("K16").Value = 2741.5715
$ws_CRP.Range("L16").Value = 7003.25
$ws_CRP.Range("M16").Value = -2454.5715
$ws_CRP.Range("N16").Value = -7577.25

$ws_CRP.Range("H31").Value = 2127.1428
$ws_CRP.Range("I31").Value = 1616.4706
$ws_CRP.Range("J31").Value = 2916.3635
$ws_CRP.Range("K31").Value = 1616.4706
$ws_CRP.Range("L31").Value = 2916.3635
$ws_CRP.Range("M31").Value = -1321.4706
$ws_CRP.Range("N31").Value = -3506.3635

$ws_CRP.Range("H34").Value = 2127.1428
$ws_CRP.Range("I34").Value = 1616.4706
$ws_CRP.Range("J34").Value = 2916.3635
$ws_CRP.Range("K34").Value = 1616.4706
$ws_CRP.Range("L34").Value = 2916.3635
$ws_CRP.Range("M34").Value = -1414.4706
$ws_CRP.Range("N34").Value = -3320.3635

$ws_CRP.Range("H99").Value = 760.1667
$ws_CRP.Range("I99").Value = 760.1667
$ws_CRP.Range("K99").Value = 760.1667
$ws_CRP.Range("M99").Value = 737.8333

$ws_CRP.Range("H113").Value = 4291.273
$ws_CRP.Range("I113").Value = 2741.5715
$ws_CRP.Range("J113").Value = 7003.25
$ws_CRP.Range("K113").Value = 2741.5715
$ws_CRP.Range("L113").Value = 7003.25
$ws_CRP.Range("M113").Value = -571.5715
$ws_CRP.Range("N113").Value = -11343.25

$ws_CRP.Range("H126").Value = 760.1667
$ws_CRP.Range("I126").Value = 760.1667
$ws_CRP.Range("K126").Value = 2280.5001
$ws_CRP.Range("M126").Value = 189.4998999999998

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H14").Value = 65.2
$ws_CUL.Range("I14").Value = 65.2
$ws_CUL.Range("K14").Value = 195.6
$ws_CUL.Range("M14").Value = -22.60000000000002

$ws_CUL.Range("H34").Value = 1677.25
$ws_CUL.Range("I34").Value = 159
$ws_CUL.Range("J34").Value = 2183.3333
$ws_CUL.Range("K34").Value = 477
$ws_CUL.Range("L34").Value = 6549.999899999999
$ws_CUL.Range("M34").Value = -393
$ws_CUL.Range("N34").Value = -6717.999899999999

$ws_CUL.Range("H35").Value = 1625
$ws_CUL.Range("J35").Value = 1625
$ws_CUL.Range("L35").Value = 4875
$ws_CUL.Range("N35").Value = -5451

$ws_CUL.Range("H39").Value = 2829.5652
$ws_CUL.Range("I39").Value = 490
$ws_CUL.Range("J39").Value = 3052.3809
$ws_CUL.Range("K39").Value = 1470
$ws_CUL.Range("L39").Value = 9157.1427
$ws_CUL.Range("M39").Value = -1176
$ws_CUL.Range("N39").Value = -9745.1427

$ws_CUL.Range("H42").Value = 1251.75
$ws_CUL.Range("I42").Value = 1001
$ws_CUL.Range("K42").Value = 3003
$ws_CUL.Range("M42").Value = -2469

$ws_CUL.Range("H55").Value = 39772.668
$ws_CUL.Range("I55").Value = 156205.2
$ws_CUL.Range("J55").Value = 3387.5
$ws_CUL.Range("K55").Value = 468615.6
$ws_CUL.Range("L55").Value = 10162.5
$ws_CUL.Range("M55").Value = -468438.6
$ws_CUL.Range("N55").Value = -10516.5

$ws_CUL.Range("H92").Value = 842
$ws_CUL.Range("I92").Value = 610.4
$ws_CUL.Range("J92").Value = 2000
$ws_CUL.Range("K92").Value = 1831.2
$ws_CUL.Range("L92").Value = 6000
$ws_CUL.Range("M92").Value = -583.1999999999998
$ws_CUL.Range("N92").Value = -8496

$ws_CUL.Range("H137").Value = 27370072
$ws_CUL.Range("I137").Value = 35715320
$ws_CUL.Range("J137").Value = 4003380
$ws_CUL.Range("K137").Value = 107145960
$ws_CUL.Range("L137").Value = 12010140
$ws_CUL.Range("M137").Value = -107140860
$ws_CUL.Range("N137").Value = -12020340

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H70").Value = 4787.8184
$ws_GSM.Range("I70").Value = 4477.722
$ws_GSM.Range("K70").Value = 4477.722
$ws_GSM.Range("M70").Value = -4207.722

$ws_GSM.Range("H73").Value = 4787.8184
$ws_GSM.Range("I73").Value = 4477.722
$ws_GSM.Range("K73").Value = 4477.722
$ws_GSM.Range("M73").Value = -3541.722

$ws_GSM.Range("H113").Value = 5681.56
$ws_GSM.Range("I113").Value = 8557.143
$ws_GSM.Range("J113").Value = 2021.7273
$ws_GSM.Range("K113").Value = 8557.143
$ws_GSM.Range("L113").Value = 2021.7273
$ws_GSM.Range("M113").Value = -6387.143
$ws_GSM.Range("N113").Value = -6361.7273

$ws_GSM.Range("H122").Value = 2754
$ws_GSM.Range("I122").Value = 2437.1428
$ws_GSM.Range("J122").Value = 3000.4443
$ws_GSM.Range("K122").Value = 7311.428400000001
$ws_GSM.Range("L122").Value = 9001.332900000001
$ws_GSM.Range("M122").Value = -4861.428400000001
$ws_GSM.Range("N122").Value = -13901.3329

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 1852.1818
$ws_LTW.Range("I40").Value = 1637.4
$ws_LTW.Range("K40").Value = 1637.4
$ws_LTW.Range("M40").Value = -1501.4

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H122").Value = 2361
$ws_WVR.Range("I122").Value = 1591.25
$ws_WVR.Range("J122").Value = 2874.1667
$ws_WVR.Range("K122").Value = 4773.75
$ws_WVR.Range("L122").Value = 8622.500100000001
$ws_WVR.Range("M122").Value = -2323.75
$ws_WVR.Range("N122").Value = -13522.5001

$ws_WVR.Range("H123").Value = 18269.75
$ws_WVR.Range("J123").Value = 18269.75
$ws_WVR.Range("L123").Value = 18269.75
$ws_WVR.Range("N123").Value = -28069.75

$ws_WVR.Range("H132").Value = 1326.34
$ws_WVR.Range("I132").Value = 1015
$ws_WVR.Range("J132").Value = 2052.8
$ws_WVR.Range("K132").Value = 3045
$ws_WVR.Range("L132").Value = 6158.400000000001
$ws_WVR.Range("M132").Value = -515
$ws_WVR.Range("N132").Value = -11218.4
